$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H (Ligand total expression value) - same correction on all rows
$ws.Range("H2").Value = 0.779989
$ws.Range("H3").Value = 0.779989
$ws.Range("H4").Value = 0.779989
$ws.Range("H5").Value = 0.779989
$ws.Range("H6").Value = 0.779989

# Row 2
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.33011
$ws.Range("N2").Value = 0.9903299999999999
$ws.Range("O2").Value = 0.04888250370346677
$ws.Range("P2").Value = 0.05014329367814145
$ws.Range("Q2").Value = 0.08582738959666665
$ws.Range("R2").Value = 0.77244650637
$ws.Range("S2").Value = 0.04888250370346677
$ws.Range("T2").Value = 0.05014329367814145

# Row 3
$ws.Range("O3").Value = 0.6674976884799527
$ws.Range("P3").Value = 0.6847139587196941
$ws.Range("S3").Value = 0.6674976884799527
$ws.Range("T3").Value = 0.6847139587196941

# Row 4
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.8214893333333334
$ws.Range("N4").Value = 2.464468
$ws.Range("O4").Value = 0.1216456798613345
$ws.Range("P4").Value = 0.1247831961915542
$ws.Range("Q4").Value = 0.2135842145391111
$ws.Range("R4").Value = 1.922257930852
$ws.Range("S4").Value = 0.1216456798613345
$ws.Range("T4").Value = 0.1247831961915542

# Row 5
$ws.Range("M5").Value = 0.509397
$ws.Range("N5").Value = 1.018794
$ws.Range("O5").Value = 0.07543122213515151
$ws.Range("P5").Value = 0.0515845089409878
$ws.Range("Q5").Value = 0.132441352211
$ws.Range("R5").Value = 0.794648113266
$ws.Range("S5").Value = 0.07543122213515151
$ws.Range("T5").Value = 0.0515845089409878

# Row 6
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5844356666666667
$ws.Range("N6").Value = 1.753307
$ws.Range("O6").Value = 0.08654290582009454
$ws.Range("P6").Value = 0.08877504246962241
$ws.Range("Q6").Value = 0.1519511304025556
$ws.Range("R6").Value = 1.367560173623
$ws.Range("S6").Value = 0.08654290582009454
$ws.Range("T6").Value = 0.08877504246962241

$wb.Save()
